$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Curve" column (P) header ---
$ws.Range("P7").Value = "Curve"

# --- Row 8 (Hamdy Abouellela): lower Final Project score, add Curve points ---
$ws.Range("N8").Value = 85
$ws.Range("P8").Value = 10

# --- Row 9 (Mahmoud Hawary): lower Final Project score, raise Penalty, add Curve points ---
$ws.Range("N9").Value = 75
$ws.Range("O9").Value = 15
$ws.Range("P9").Value = 10

# --- Row 10 (Khadija Mahbuba): lower Final Project score, add Curve points ---
$ws.Range("N10").Value = 85
$ws.Range("P10").Value = 10

# --- Row 12: recompute helper row for the new Curve column ---
$ws.Range("O12").Formula = "=N12-AVERAGE(O8:O10)"
$ws.Range("P12").Formula = "=O12+AVERAGE(P8:P10)"

# --- Row 13: weight for new Curve column ---
$ws.Range("P13").Value = 20

# --- Row 16: weight-check row - move the 20 from Penalty(O) to the new Curve(P) column ---
$ws.Range("O16").ClearContents()
$ws.Range("P16").Value = 20
$ws.Range("Q16").Formula = "=SUM(D16:P16)"

# --- Row 17: replace Penalty-only adjustment with Penalty+Curve adjustment ---
$ws.Range("O17").ClearContents()
$ws.Range("P17").Formula = "=(N8-O8+P8)*P`$13/100"
$ws.Range("Q17").Formula = "=SUM(D17:P17)/`$Q`$16"

# --- Row 18: same change for Mahmoud Hawary ---
$ws.Range("O18").ClearContents()
$ws.Range("P18").Formula = "=(N9-O9+P9)*P`$13/100"
$ws.Range("Q18").Formula = "=SUM(D18:P18)/`$Q`$16"

# --- Row 19: same change for Khadija Mahbuba ---
$ws.Range("O19").ClearContents()
$ws.Range("P19").Formula = "=(N10-O10+P10)*P`$13/100"
$ws.Range("Q19").Formula = "=SUM(D19:P19)/`$Q`$16"

# --- Formatting: the Final Project scores are no longer flagged red, ---
# --- while the grade-scale notes next to the curved rows are now flagged red instead ---
$ws.Range("N8").ClearFormats()
$ws.Range("N9").ClearFormats()
$ws.Range("N10").ClearFormats()

$ws.Range("Q9").ClearFormats()
$ws.Range("Q9").Font.Color = 255

$ws.Range("Q10").ClearFormats()
$ws.Range("Q10").Font.Color = 255

# --- Move the active selection like the author left it ---
$ws.Range("P21").Select()
